$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2118.7742
$ws.Range("I40").Value = 1823.6471
$ws.Range("J40").Value = 2477.1428
$ws.Range("K40").Value = 1823.6471
$ws.Range("L40").Value = 2477.1428
$ws.Range("M40").Value = -1648.6471
$ws.Range("N40").Value = -2827.1428
$ws.Range("H48").Value = 4000
$ws.Range("J48").Value = 4000
$ws.Range("L48").Value = 12000
$ws.Range("N48").Value = -12584
$ws.Range("H56").Value = 4000
$ws.Range("J56").Value = 4000
$ws.Range("L56").Value = 12000
$ws.Range("N56").Value = -13068
$ws.Range("H80").Value = 1106.6666
$ws.Range("I80").Value = 2090
$ws.Range("J80").Value = 910
$ws.Range("K80").Value = 6270
$ws.Range("L80").Value = 2730
$ws.Range("M80").Value = -5272
$ws.Range("N80").Value = -4726
$ws.Range("H83").Value = 1106.6666
$ws.Range("I83").Value = 2090
$ws.Range("J83").Value = 910
$ws.Range("K83").Value = 18810
$ws.Range("L83").Value = 8190
$ws.Range("M83").Value = -13818
$ws.Range("N83").Value = -18174
$ws.Range("H112").Value = 2361.5
$ws.Range("J112").Value = 3034.75
$ws.Range("L112").Value = 9104.25
$ws.Range("N112").Value = -11320.25
$ws.Range("H116").Value = 3214.8333
$ws.Range("I116").Value = 2095.3333
$ws.Range("J116").Value = 4334.3335
$ws.Range("K116").Value = 2095.3333
$ws.Range("L116").Value = 4334.3335
$ws.Range("M116").Value = 1346.6667
$ws.Range("N116").Value = -11218.3335
$ws.Range("H121").Value = 599.61536
$ws.Range("J121").Value = 599.61536
$ws.Range("L121").Value = 1798.84608
$ws.Range("N121").Value = -5292.84608
$ws.Range("H129").Value = 824.2432
$ws.Range("J129").Value = 903.4516
$ws.Range("L129").Value = 2710.3548
$ws.Range("N129").Value = -12710.3548
$ws.Range("H137").Value = 1531.5116
$ws.Range("I137").Value = 1332.1034
$ws.Range("J137").Value = 1944.5714
$ws.Range("K137").Value = 3996.3102
$ws.Range("L137").Value = 5833.7142
$ws.Range("M137").Value = -1446.3102
$ws.Range("N137").Value = -10933.7142

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5710
$ws.Range("I32").Value = 5048.3506
$ws.Range("J32").Value = 9349.071
$ws.Range("K32").Value = 5048.3506
$ws.Range("L32").Value = 9349.071
$ws.Range("M32").Value = -4761.3506
$ws.Range("N32").Value = -9923.071
$ws.Range("H74").Value = 3694
$ws.Range("I74").Value = 2803
$ws.Range("J74").Value = 4288
$ws.Range("K74").Value = 2803
$ws.Range("L74").Value = 4288
$ws.Range("M74").Value = -1929
$ws.Range("N74").Value = -6036
$ws.Range("H77").Value = 3694
$ws.Range("I77").Value = 2803
$ws.Range("J77").Value = 4288
$ws.Range("K77").Value = 14015
$ws.Range("L77").Value = 21440
$ws.Range("M77").Value = -9647
$ws.Range("N77").Value = -30176
$ws.Range("H125").Value = 35707.5
$ws.Range("J125").Value = 35707.5
$ws.Range("L125").Value = 35707.5
$ws.Range("N125").Value = -45547.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 350
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = -27

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1668.7567
$ws.Range("I31").Value = 1576.2222
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1576.2222
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1281.2222
$ws.Range("N31").Value = -5590
$ws.Range("H34").Value = 1668.7567
$ws.Range("I34").Value = 1576.2222
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1576.2222
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1374.2222
$ws.Range("N34").Value = -5404
$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15368
$ws.Range("H101").Value = 15000
$ws.Range("J101").Value = 15000
$ws.Range("L101").Value = 15000
$ws.Range("N101").Value = -21490
$ws.Range("H105").Value = 903
$ws.Range("I105").Value = 885
$ws.Range("J105").Value = 1011
$ws.Range("K105").Value = 885
$ws.Range("L105").Value = 1011
$ws.Range("M105").Value = 862
$ws.Range("N105").Value = -4505
$ws.Range("H131").Value = 15279.777
$ws.Range("J131").Value = 20771.666
$ws.Range("L131").Value = 20771.666
$ws.Range("N131").Value = -30851.666

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10485.177
$ws.Range("I3").Value = 6208.1816
$ws.Range("J3").Value = 18326.334
$ws.Range("K3").Value = 18624.5448
$ws.Range("L3").Value = 54979.00199999999
$ws.Range("M3").Value = -18512.5448
$ws.Range("N3").Value = -55203.00199999999
$ws.Range("H54").Value = 1500
$ws.Range("J54").Value = 1500
$ws.Range("L54").Value = 4500
$ws.Range("N54").Value = -5618
$ws.Range("H131").Value = 14707230
$ws.Range("I131").Value = 125000500
$ws.Range("J131").Value = 1460.9667
$ws.Range("K131").Value = 375001500
$ws.Range("L131").Value = 4382.9001
$ws.Range("M131").Value = -374996460
$ws.Range("N131").Value = -14462.9001

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5110
$ws.Range("I80").Value = 4385.7144
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 4385.7144
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -3387.7144
$ws.Range("N80").Value = -8796
$ws.Range("H83").Value = 5110
$ws.Range("I83").Value = 4385.7144
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 21928.572
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -16936.572
$ws.Range("N83").Value = -43984
$ws.Range("H105").Value = 26000
$ws.Range("J105").Value = 26000
$ws.Range("L105").Value = 26000
$ws.Range("N105").Value = -32988
$ws.Range("H123").Value = 21444.428
$ws.Range("J123").Value = 21444.428
$ws.Range("L123").Value = 21444.428
$ws.Range("N123").Value = -26344.428
$ws.Range("H132").Value = 3274.9583
$ws.Range("I132").Value = 3117.0625
$ws.Range("K132").Value = 9351.1875
$ws.Range("M132").Value = -6821.1875

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7489.5
$ws.Range("I40").Value = 3788.6667
$ws.Range("J40").Value = 9710
$ws.Range("K40").Value = 3788.6667
$ws.Range("L40").Value = 9710
$ws.Range("M40").Value = -3652.6667
$ws.Range("N40").Value = -9982
$ws.Range("H55").Value = 166.19444
$ws.Range("I55").Value = 103.78125
$ws.Range("J55").Value = 665.5
$ws.Range("K55").Value = 103.78125
$ws.Range("L55").Value = 665.5
$ws.Range("M55").Value = 69.21875
$ws.Range("N55").Value = -1011.5
$ws.Range("H101").Value = 13975
$ws.Range("J101").Value = 13975
$ws.Range("L101").Value = 13975
$ws.Range("N101").Value = -20465
$ws.Range("H106").Value = 22637.666
$ws.Range("J106").Value = 22637.666
$ws.Range("L106").Value = 22637.666
$ws.Range("N106").Value = -25161.666
$ws.Range("H136").Value = 1475.0834
$ws.Range("I136").Value = 1343.0476
$ws.Range("J136").Value = 2399.3333
$ws.Range("K136").Value = 4029.142800000001
$ws.Range("L136").Value = 7197.999899999999
$ws.Range("M136").Value = -1479.142800000001
$ws.Range("N136").Value = -12297.9999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22729272
$ws.Range("I122").Value = 22729272
$ws.Range("K122").Value = 68187816
$ws.Range("M122").Value = -68185366
$ws.Range("H126").Value = 71429030
$ws.Range("I126").Value = 71429030
$ws.Range("K126").Value = 214287090
$ws.Range("M126").Value = -214284620
